# WR_89787325_WeekEnding_071325.xlsx
# -----------------------------------------------------------------------------
# Business fix: a single Work Request Excel export must not contain a day
# section (Friday 07/11/2025) that only carries zeroed-out placeholder rows.
# Remove that Friday block entirely (header + column headers + its one data
# row + its TOTAL row) so the report starts with the first day that actually
# has billed units (Saturday 07/12/2025), followed by Sunday (07/13/2025).
# Deleting the whole rows lets Excel shift the remaining sections (and their
# existing per-row styles / merged cells) up automatically.
# -----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the entire "Friday (07/11/2025)" section: row 14 (day header),
#    row 15 (column headers), row 16 (the single "Point 08" line with all
#    zero units/pricing), row 17 (its TOTAL row) and the two blank spacer
#    rows 18-19 that preceded the next section. This shifts the Saturday
#    section up to start at row 14 and the Sunday section up to row 23.
$ws.Rows("14:19").Delete()

# 2) Report header refresh.
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:00 AM"

# 3) Report summary numbers now that Friday's zero-value line item is gone.
$ws.Range("C8").Value = 6485.299999999999
$ws.Range("C9").Value = 6

# 4) Scope ID is no longer known for this export - clear it.
$ws.Range("G10").Value = ""

# 5) Saturday (07/12/2025) section now occupies rows 14-20; fill in the real
#    computed pricing for each line (was 0 while Friday's block was still
#    first) and the section TOTAL.
$ws.Range("H16").Value = 1297.06   # Point 22
$ws.Range("H17").Value = 1297.06   # Point 20
$ws.Range("H18").Value = 1297.06   # Point 18
$ws.Range("H19").Value = 648.53    # Point 16
$ws.Range("H20").Value = 4539.71   # Saturday TOTAL

# 6) Sunday (07/13/2025) section now occupies rows 23-27; fill in the real
#    computed pricing for each line and the section TOTAL.
$ws.Range("H25").Value = 1297.06   # Point 14
$ws.Range("H26").Value = 648.53    # Point 12
$ws.Range("H27").Value = 1945.59   # Sunday TOTAL
